$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unified_table")

# Correct the hardcoded row-count value in A43 (was 22, should be 23).
# All following rows (44:74) derive their value from formulas chained
# off A43, so fixing this one cell cascades the correct sequential
# numbering through the rest of the "Unified_table" sheet.
$ws.Range("A43").Value = 23

$excel.CalculateFullRebuild()
